# "Natmi following Dr Hou advice"
# The Alcam-Nrp1 NATMI ligand-receptor results now include "ECs" as a third
# sending cluster (previously only FAPs/sCs sent to FAPs/sCs/ECs). That grows
# the table from 2 sending clusters x 3 target clusters (6 data rows) to
# 3 x 3 (9 data rows), so every metric column is recomputed and three new
# rows are appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2-10 (A:T): Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric NATMI metric columns.
$rowData = @{
    2 = @("ECs", "Alcam", "Nrp1", "ECs", 1, 0.3333333333333333, 62.12558000000001, 186.37674, 0.9736910227596813, 0.9736910227596813, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 6388.307807496049, 57494.77026746443, 0.5412865624319931, 0.541286562431993)
    3 = @("ECs", "Alcam", "Nrp1", "FAPs", 1, 0.3333333333333333, 62.12558000000001, 186.37674, 0.9736910227596813, 0.9736910227596813, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 3955.077709532767, 35595.69938579491, 0.335116979653413, 0.335116979653413)
    4 = @("ECs", "Alcam", "Nrp1", "sCs", 1, 0.3333333333333333, 62.12558000000001, 186.37674, 0.9736910227596813, 0.9736910227596813, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 1148.1947188393, 10333.7524695537, 0.09728748067427527, 0.09728748067427526)
    5 = @("FAPs", "Alcam", "Nrp1", "ECs", 3, 1, 0.5683613333333334, 1.705084, 0.008907897969731461, 0.008907897969731461, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 58.44399590655246, 525.9959631589721, 0.00495200772917153, 0.004952007729171528)
    6 = @("FAPs", "Alcam", "Nrp1", "FAPs", 3, 1, 0.5683613333333334, 1.705084, 0.008907897969731461, 0.008907897969731461, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 36.18337632303778, 325.65038690734, 0.00306584716598949, 0.00306584716598949)
    7 = @("FAPs", "Alcam", "Nrp1", "sCs", 3, 1, 0.5683613333333334, 1.705084, 0.008907897969731461, 0.008907897969731461, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 10.50436038304667, 94.53924344741999, 0.0008900430745704425, 0.0008900430745704424)
    8 = @("sCs", "Alcam", "Nrp1", "ECs", 3, 1, 1.110262, 3.330786, 0.01740107927058724, 0.01740107927058724, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 114.1670693934154, 1027.503624540738, 0.009673469469079718, 0.009673469469079718)
    9 = @("sCs", "Alcam", "Nrp1", "FAPs", 3, 1, 1.110262, 3.330786, 0.01740107927058724, 0.01740107927058724, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 70.68219705862334, 636.13977352761, 0.005988960554798161, 0.005988960554798163)
    10 = @("sCs", "Alcam", "Nrp1", "sCs", 3, 1, 1.110262, 3.330786, 0.01740107927058724, 0.01740107927058724, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 20.51967909077, 184.67711181693, 0.001738649246709362, 0.001738649246709362)
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rowData.Keys | Sort-Object) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range(("{0}{1}" -f $cols[$i], $r)).Value = $vals[$i]
    }
}
